$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the disclaimer date text (shared string) from 2021-05-27 to 2021-05-28
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."

# Update the Weight / Percent Change numeric values
$ws.Range("D2").Value = 0.2476830222942763
$ws.Range("E2").Value = 0.00148001973359646

$ws.Range("D3").Value = 0.4969079969568
$ws.Range("E3").Value = 0.002359108781127217

$ws.Range("D4").Value = 0.09684787783769284
$ws.Range("E4").Value = 0.001945903872348609

$ws.Range("D5").Value = 0.1016374206797095
$ws.Range("E5").Value = 0.0009020541988731257

$ws.Range("D6").Value = 0.05692368223152135
$ws.Range("E6").Value = -0.000223164472216264

$ws.Range("E7").Value = 0.001806271558803418
